$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '22.366.45'

# Row 3
$ws.Range("D3").Value = '1.567.07'
$ws.Range("E3").Value = '  -0.06%  '

# Row 4
$ws.Range("E4").Value = '  +0.39%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '1.003'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.28%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '290.04'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.28%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.3745'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.99%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '49.08'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -0.12%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.3373'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -0.87%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.07515'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -2.00%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '1.127'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -3.78%  '

# Row 12
$ws.Range("E12").Value = '  +0.43%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '20.80'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -3.61%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.896'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -2.66%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '6.866'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -1.04%  '

# Row 16
$ws.Range("D16").Value = '1.564.98'
$ws.Range("E16").Value = '  +0.08%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.00001116'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -1.37%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '89.43'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -0.89%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06724'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.31%  '

# Row 20
$ws.Range("E20").Value = '  +0.40%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.172'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -1.35%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '16.38'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -1.15%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '11.83'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -1.81%  '

# Row 24
$ws.Range("D24").Value = '22.371.02'
$ws.Range("E24").Value = '  -0.11%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.379'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.43%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.712'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -3.46%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '147.61'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +1.44%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '4.990'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +0.26%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '124.84'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -0.61%  '

# Row 31
$ws.Range("D31").Value = '1.738.88'
$ws.Range("E31").Value = '  -0.12%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '2.015'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +0.28%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.9813'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -3.31%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.946'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -4.63%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '9.918'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -1.74%  '

# Row 36
$ws.Range("B36").Value = 'TrustWalletToken'
$ws.Range("C36").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.404'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +9.33%  '

# Row 37
$ws.Range("B37").Value = 'Stellar'
$ws.Range("C37").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.08419'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.92%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.02454'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -3.49%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.2271'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -2.50%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.06401'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -0.54%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '5.359'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -3.53%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.6240'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -1.90%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '10.92'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -7.02%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.004'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.36%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '13.86'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -2.51%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.791'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.66%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.5832'
$ws.Range("D47").Style = "Normal"

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.046'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -2.73%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.247'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -1.67%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '124.19'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -0.35%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.07309'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +0.26%  '
